$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 190.44444
$ws.Range("I2").Value = 202.8
$ws.Range("J2").Value = 175
$ws.Range("K2").Value = 202.8
$ws.Range("L2").Value = 175
$ws.Range("M2").Value = -89.80000000000001
$ws.Range("N2").Value = -401
$ws.Range("H33").Value = 15626627
$ws.Range("I33").Value = 27779116
$ws.Range("J33").Value = 1997.8572
$ws.Range("K33").Value = 27779116
$ws.Range("L33").Value = 1997.8572
$ws.Range("M33").Value = -27778887
$ws.Range("N33").Value = -2455.8572
$ws.Range("H40").Value = 8133.6665
$ws.Range("J40").Value = 11200.5
$ws.Range("L40").Value = 11200.5
$ws.Range("N40").Value = -11550.5
$ws.Range("H43").Value = 13624
$ws.Range("I43").Value = 9000
$ws.Range("K43").Value = 9000
$ws.Range("M43").Value = -8931
$ws.Range("H45").Value = 2800
$ws.Range("J45").Value = 2800
$ws.Range("L45").Value = 8400
$ws.Range("N45").Value = -8784
$ws.Range("H98").Value = 1961.8182
$ws.Range("I98").Value = 1698.7368
$ws.Range("J98").Value = 3628
$ws.Range("K98").Value = 1698.7368
$ws.Range("L98").Value = 3628
$ws.Range("M98").Value = -200.7367999999999
$ws.Range("N98").Value = -6624
$ws.Range("H101").Value = 1404.5834
$ws.Range("J101").Value = 878.3333
$ws.Range("L101").Value = 2634.9999
$ws.Range("N101").Value = -5878.9999
$ws.Range("H113").Value = 21146.8
$ws.Range("I113").Value = 21050.691
$ws.Range("K113").Value = 21050.691
$ws.Range("M113").Value = -17796.691
$ws.Range("H116").Value = 3848.4312
$ws.Range("I116").Value = 3579.4119
$ws.Range("J116").Value = 5808.4287
$ws.Range("K116").Value = 3579.4119
$ws.Range("L116").Value = 5808.4287
$ws.Range("M116").Value = -137.4119000000001
$ws.Range("N116").Value = -12692.4287
$ws.Range("H122").Value = 1961.8182
$ws.Range("I122").Value = 1698.7368
$ws.Range("J122").Value = 3628
$ws.Range("K122").Value = 5096.2104
$ws.Range("L122").Value = 10884
$ws.Range("M122").Value = -2646.2104
$ws.Range("N122").Value = -15784
$ws.Range("H132").Value = 5710.769
$ws.Range("I132").Value = 4291.281
$ws.Range("J132").Value = 15824.625
$ws.Range("K132").Value = 12873.843
$ws.Range("L132").Value = 47473.875
$ws.Range("M132").Value = -10343.843
$ws.Range("N132").Value = -52533.875
$ws.Range("H137").Value = 11149.167
$ws.Range("I137").Value = 4862.9585
$ws.Range("J137").Value = 16178.134
$ws.Range("K137").Value = 14588.8755
$ws.Range("L137").Value = 48534.402
$ws.Range("M137").Value = -12038.8755
$ws.Range("N137").Value = -53634.402
$ws.Range("H138").Value = 6017.09
$ws.Range("I138").Value = 2384.3713
$ws.Range("J138").Value = 7973.1694
$ws.Range("K138").Value = 7153.113899999999
$ws.Range("L138").Value = 23919.5082
$ws.Range("M138").Value = -2013.113899999999
$ws.Range("N138").Value = -34199.5082

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1180004.9
$ws.Range("I32").Value = 1858046.4
$ws.Range("K32").Value = 1858046.4
$ws.Range("M32").Value = -1857759.4
$ws.Range("H45").Value = 3169.2812
$ws.Range("I45").Value = 3363.2693
$ws.Range("K45").Value = 3363.2693
$ws.Range("M45").Value = -2986.2693
$ws.Range("H61").Value = 17324.309
$ws.Range("I61").Value = 7241.4585
$ws.Range("J61").Value = 33456.867
$ws.Range("K61").Value = 7241.4585
$ws.Range("L61").Value = 33456.867
$ws.Range("M61").Value = -7029.4585
$ws.Range("N61").Value = -33880.867
$ws.Range("H63").Value = 2793.8
$ws.Range("I63").Value = 1284.5
$ws.Range("J63").Value = 3800
$ws.Range("K63").Value = 1284.5
$ws.Range("L63").Value = 3800
$ws.Range("M63").Value = -598.5
$ws.Range("N63").Value = -5172
$ws.Range("H66").Value = 2793.8
$ws.Range("I66").Value = 1284.5
$ws.Range("J66").Value = 3800
$ws.Range("K66").Value = 6422.5
$ws.Range("L66").Value = 19000
$ws.Range("M66").Value = -2990.5
$ws.Range("N66").Value = -25864
$ws.Range("H74").Value = 20683.75
$ws.Range("I74").Value = 3878.4
$ws.Range("K74").Value = 3878.4
$ws.Range("M74").Value = -3004.4
$ws.Range("H77").Value = 20683.75
$ws.Range("I77").Value = 3878.4
$ws.Range("K77").Value = 19392
$ws.Range("M77").Value = -15024
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = None
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = None
$ws.Range("N83").ClearContents()
$ws.Range("I97").Value = 2550
$ws.Range("J97").Value = 1828.5714
$ws.Range("K97").Value = 2550
$ws.Range("L97").Value = 1828.5714
$ws.Range("M97").Value = -2054
$ws.Range("N97").Value = -2820.5714
$ws.Range("H122").Value = 4236.9473
$ws.Range("I122").Value = 2777.5
$ws.Range("J122").Value = 7399.0835
$ws.Range("K122").Value = 8332.5
$ws.Range("L122").Value = 22197.2505
$ws.Range("M122").Value = -5882.5
$ws.Range("N122").Value = -27097.2505
$ws.Range("H132").Value = 10020.658
$ws.Range("I132").Value = 4698.9653
$ws.Range("J132").Value = 22881.416
$ws.Range("K132").Value = 14096.8959
$ws.Range("L132").Value = 68644.24800000001
$ws.Range("M132").Value = -11566.8959
$ws.Range("N132").Value = -73704.24800000001
$ws.Range("H136").Value = 17324.309
$ws.Range("I136").Value = 7241.4585
$ws.Range("J136").Value = 33456.867
$ws.Range("K136").Value = 21724.3755
$ws.Range("L136").Value = 100370.601
$ws.Range("M136").Value = -19174.3755
$ws.Range("N136").Value = -105470.601

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 17069.166
$ws.Range("I20").Value = 4248.6924
$ws.Range("K20").Value = 4248.6924
$ws.Range("M20").Value = -4001.6924
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H86").Value = 4486.2173
$ws.Range("I86").Value = 3688.5881
$ws.Range("K86").Value = 3688.5881
$ws.Range("M86").Value = -2565.5881
$ws.Range("H89").Value = 4486.2173
$ws.Range("I89").Value = 3688.5881
$ws.Range("K89").Value = 18442.9405
$ws.Range("M89").Value = -12826.9405
$ws.Range("H134").Value = 6883.3276
$ws.Range("I134").Value = 2229.4211
$ws.Range("J134").Value = 15725.75
$ws.Range("K134").Value = 6688.263300000001
$ws.Range("L134").Value = 47177.25
$ws.Range("M134").Value = -4153.263300000001
$ws.Range("N134").Value = -52247.25
$ws.Range("H141").Value = 49666.668
$ws.Range("J141").Value = 49666.668
$ws.Range("L141").Value = 49666.668
$ws.Range("N141").Value = -60026.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 358597.5
$ws.Range("I6").Value = 714669.3
$ws.Range("K6").Value = 714669.3
$ws.Range("M6").Value = -714556.3
$ws.Range("H7").Value = 626.7568
$ws.Range("I7").Value = 501.32257
$ws.Range("J7").Value = 1274.8334
$ws.Range("K7").Value = 501.32257
$ws.Range("L7").Value = 1274.8334
$ws.Range("M7").Value = -388.32257
$ws.Range("N7").Value = -1500.8334
$ws.Range("H16").Value = 9625.546
$ws.Range("I16").Value = 3769.5
$ws.Range("K16").Value = 3769.5
$ws.Range("M16").Value = -3482.5
$ws.Range("H17").Value = 899.5
$ws.Range("I17").Value = 899.5
$ws.Range("K17").Value = 899.5
$ws.Range("M17").Value = -725.5
$ws.Range("H22").Value = 3766.3635
$ws.Range("I22").Value = 3363.8333
$ws.Range("J22").Value = 4249.4
$ws.Range("K22").Value = 3363.8333
$ws.Range("L22").Value = 4249.4
$ws.Range("M22").Value = -3013.8333
$ws.Range("N22").Value = -4949.4
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("H31").Value = 20761.834
$ws.Range("I31").Value = 6132.731
$ws.Range("J31").Value = 58797.5
$ws.Range("K31").Value = 6132.731
$ws.Range("L31").Value = 58797.5
$ws.Range("M31").Value = -5837.731
$ws.Range("N31").Value = -59387.5
$ws.Range("H34").Value = 20761.834
$ws.Range("I34").Value = 6132.731
$ws.Range("J34").Value = 58797.5
$ws.Range("K34").Value = 6132.731
$ws.Range("L34").Value = 58797.5
$ws.Range("M34").Value = -5930.731
$ws.Range("N34").Value = -59201.5
$ws.Range("H51").Value = 35000
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472
$ws.Range("H56").Value = 13599.6
$ws.Range("I56").Value = 5999.6665
$ws.Range("J56").Value = 24999.5
$ws.Range("K56").Value = 5999.6665
$ws.Range("L56").Value = 24999.5
$ws.Range("M56").Value = -5154.6665
$ws.Range("N56").Value = -26689.5
$ws.Range("H58").Value = 13084.703
$ws.Range("I58").Value = 5139.8696
$ws.Range("J58").Value = 26136.928
$ws.Range("K58").Value = 5139.8696
$ws.Range("L58").Value = 26136.928
$ws.Range("M58").Value = -4936.8696
$ws.Range("N58").Value = -26542.928
$ws.Range("H61").Value = 35000
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696
$ws.Range("H62").Value = 6410.125
$ws.Range("J62").Value = 8003
$ws.Range("L62").Value = 8003
$ws.Range("N62").Value = -9251
$ws.Range("H65").Value = 6410.125
$ws.Range("J65").Value = 8003
$ws.Range("L65").Value = 40015
$ws.Range("N65").Value = -46255
$ws.Range("H68").Value = 49750
$ws.Range("I68").Value = 49750
$ws.Range("K68").Value = 49750
$ws.Range("M68").Value = -49001
$ws.Range("H70").Value = 17000
$ws.Range("J70").Value = 17000
$ws.Range("L70").Value = 17000
$ws.Range("N70").Value = -17630
$ws.Range("H71").Value = 49750
$ws.Range("I71").Value = 49750
$ws.Range("K71").Value = 149250
$ws.Range("M71").Value = -145506
$ws.Range("H73").Value = 17000
$ws.Range("J73").Value = 17000
$ws.Range("L73").Value = 17000
$ws.Range("N73").Value = -19184
$ws.Range("H99").Value = 11119.728
$ws.Range("I99").Value = 5329.5
$ws.Range("J99").Value = 14428.429
$ws.Range("K99").Value = 5329.5
$ws.Range("L99").Value = 14428.429
$ws.Range("M99").Value = -3831.5
$ws.Range("N99").Value = -17424.429
$ws.Range("H113").Value = 9625.546
$ws.Range("I113").Value = 3769.5
$ws.Range("K113").Value = 3769.5
$ws.Range("M113").Value = -1599.5
$ws.Range("H126").Value = 11119.728
$ws.Range("I126").Value = 5329.5
$ws.Range("J126").Value = 14428.429
$ws.Range("K126").Value = 15988.5
$ws.Range("L126").Value = 43285.287
$ws.Range("M126").Value = -13518.5
$ws.Range("N126").Value = -48225.287
$ws.Range("H132").Value = 5815.0205
$ws.Range("I132").Value = 1905.6428
$ws.Range("J132").Value = 11027.523
$ws.Range("K132").Value = 5716.928400000001
$ws.Range("L132").Value = 33082.569
$ws.Range("M132").Value = -3186.928400000001
$ws.Range("N132").Value = -38142.569
$ws.Range("H136").Value = 13084.703
$ws.Range("I136").Value = 5139.8696
$ws.Range("J136").Value = 26136.928
$ws.Range("K136").Value = 15419.6088
$ws.Range("L136").Value = 78410.784
$ws.Range("M136").Value = -12869.6088
$ws.Range("N136").Value = -83510.784

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 2194719.5
$ws.Range("I50").Value = 1985.2858
$ws.Range("K50").Value = 5955.857400000001
$ws.Range("M50").Value = -5474.857400000001
$ws.Range("H53").Value = 2194719.5
$ws.Range("I53").Value = 1985.2858
$ws.Range("K53").Value = 5955.857400000001
$ws.Range("M53").Value = -5474.857400000001
$ws.Range("H113").Value = 2632.6667
$ws.Range("J113").Value = 2699
$ws.Range("L113").Value = 8097
$ws.Range("N113").Value = -12437
$ws.Range("H114").Value = 3288.125
$ws.Range("J114").Value = 6166.3335
$ws.Range("L114").Value = 18499.0005
$ws.Range("N114").Value = -25007.0005
$ws.Range("H117").Value = 25800.334
$ws.Range("J117").Value = 38543
$ws.Range("L117").Value = 115629
$ws.Range("N117").Value = -122513
$ws.Range("H121").Value = 10774.357
$ws.Range("J121").Value = 16454.666
$ws.Range("L121").Value = 49363.99800000001
$ws.Range("N121").Value = -51983.99800000001
$ws.Range("H122").Value = 15374118
$ws.Range("J122").Value = 2836775.8
$ws.Range("L122").Value = 25530982.2
$ws.Range("N122").Value = -25535882.2
$ws.Range("H127").Value = 71956.7
$ws.Range("J127").Value = 71956.7
$ws.Range("L127").Value = 215870.1
$ws.Range("N127").Value = -225790.1
$ws.Range("H129").Value = 7695438
$ws.Range("I129").Value = 25000820
$ws.Range("J129").Value = 4156.8887
$ws.Range("K129").Value = 75002460
$ws.Range("L129").Value = 12470.6661
$ws.Range("M129").Value = -74997460
$ws.Range("N129").Value = -22470.6661
$ws.Range("H131").Value = 1466.9
$ws.Range("I131").Value = 800
$ws.Range("J131").Value = 1480.5103
$ws.Range("K131").Value = 2400
$ws.Range("L131").Value = 4441.5309
$ws.Range("M131").Value = 2640
$ws.Range("N131").Value = -14521.5309
$ws.Range("H134").Value = 5042.9604
$ws.Range("I134").Value = 1051.6
$ws.Range("K134").Value = 3154.8
$ws.Range("M134").Value = 1915.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3341.1667
$ws.Range("I2").Value = 4922.5
$ws.Range("J2").Value = 178.5
$ws.Range("K2").Value = 4922.5
$ws.Range("L2").Value = 178.5
$ws.Range("M2").Value = -4809.5
$ws.Range("N2").Value = -404.5
$ws.Range("H15").Value = 8580
$ws.Range("J15").Value = 8598
$ws.Range("L15").Value = 8598
$ws.Range("N15").Value = -9174
$ws.Range("H81").Value = 8580
$ws.Range("J81").Value = 8598
$ws.Range("L81").Value = 8598
$ws.Range("N81").Value = -10594
$ws.Range("H84").Value = 8580
$ws.Range("J84").Value = 8598
$ws.Range("L84").Value = 25794
$ws.Range("N84").Value = -35778
$ws.Range("H97").Value = 693.5909
$ws.Range("I97").Value = 676.4737
$ws.Range("K97").Value = 676.4737
$ws.Range("M97").Value = -180.4737
$ws.Range("H132").Value = 7135.5864
$ws.Range("I132").Value = 3392.76
$ws.Range("J132").Value = 30528.25
$ws.Range("K132").Value = 10178.28
$ws.Range("L132").Value = 91584.75
$ws.Range("M132").Value = -7648.280000000001
$ws.Range("N132").Value = -96644.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 2722.111
$ws.Range("I3").Value = 1999.8334
$ws.Range("J3").Value = 4166.6665
$ws.Range("K3").Value = 1999.8334
$ws.Range("L3").Value = 4166.6665
$ws.Range("M3").Value = -1887.8334
$ws.Range("N3").Value = -4390.6665
$ws.Range("H15").Value = 2722.111
$ws.Range("I15").Value = 1999.8334
$ws.Range("J15").Value = 4166.6665
$ws.Range("K15").Value = 1999.8334
$ws.Range("L15").Value = 4166.6665
$ws.Range("M15").Value = -1829.8334
$ws.Range("N15").Value = -4506.6665
$ws.Range("H22").Value = 8630.714
$ws.Range("I22").Value = 8589.8125
$ws.Range("J22").Value = 8685.25
$ws.Range("K22").Value = 8589.8125
$ws.Range("L22").Value = 8685.25
$ws.Range("M22").Value = -8294.8125
$ws.Range("N22").Value = -9275.25
$ws.Range("H27").Value = 8630.714
$ws.Range("I27").Value = 8589.8125
$ws.Range("J27").Value = 8685.25
$ws.Range("K27").Value = 8589.8125
$ws.Range("L27").Value = 8685.25
$ws.Range("M27").Value = -8482.8125
$ws.Range("N27").Value = -8899.25
$ws.Range("H42").Value = 14800
$ws.Range("I42").Value = 14800
$ws.Range("K42").Value = 14800
$ws.Range("M42").Value = -14237
$ws.Range("H46").Value = 4232.75
$ws.Range("I46").Value = 4900
$ws.Range("K46").Value = 4900
$ws.Range("M46").Value = -4712
$ws.Range("H49").Value = 14800
$ws.Range("I49").Value = 14800
$ws.Range("K49").Value = 14800
$ws.Range("M49").Value = -14653
$ws.Range("H68").Value = 3076.3845
$ws.Range("I68").Value = 2798.6667
$ws.Range("J68").Value = 3314.4285
$ws.Range("K68").Value = 2798.6667
$ws.Range("L68").Value = 3314.4285
$ws.Range("M68").Value = -2049.6667
$ws.Range("N68").Value = -4812.4285
$ws.Range("H71").Value = 3076.3845
$ws.Range("I71").Value = 2798.6667
$ws.Range("J71").Value = 3314.4285
$ws.Range("K71").Value = 13993.3335
$ws.Range("L71").Value = 16572.1425
$ws.Range("M71").Value = -10249.3335
$ws.Range("N71").Value = -24060.1425
$ws.Range("H100").Value = 3849.75
$ws.Range("I100").Value = 3200
$ws.Range("J100").Value = 4499.5
$ws.Range("K100").Value = 3200
$ws.Range("L100").Value = 4499.5
$ws.Range("M100").Value = -2659
$ws.Range("N100").Value = -5581.5
$ws.Range("H122").Value = 8190.231
$ws.Range("I122").Value = 6412.5
$ws.Range("K122").Value = 19237.5
$ws.Range("M122").Value = -16787.5
$ws.Range("H132").Value = 7873.925
$ws.Range("I132").Value = 4575.815
$ws.Range("J132").Value = 14723.846
$ws.Range("K132").Value = 13727.445
$ws.Range("L132").Value = 44171.538
$ws.Range("M132").Value = -11197.445
$ws.Range("N132").Value = -49231.538

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 1800
$ws.Range("H68").Value = 48333.332
$ws.Range("J68").Value = 48333.332
$ws.Range("L68").Value = 48333.332
$ws.Range("N68").Value = -49955.332
$ws.Range("H69").Value = 64999.2
$ws.Range("J69").Value = 64999.2
$ws.Range("L69").Value = 64999.2
$ws.Range("N69").Value = -66497.2
$ws.Range("H71").Value = 48333.332
$ws.Range("J71").Value = 48333.332
$ws.Range("L71").Value = 144999.996
$ws.Range("N71").Value = -153111.996
$ws.Range("H72").Value = 64999.2
$ws.Range("J72").Value = 64999.2
$ws.Range("L72").Value = 194997.6
$ws.Range("N72").Value = -202485.6
$ws.Range("H96").Value = 2968.3
$ws.Range("I96").Value = 3098.75
$ws.Range("J96").Value = 2881.3333
$ws.Range("K96").Value = 3098.75
$ws.Range("L96").Value = 2881.3333
$ws.Range("M96").Value = -1725.75
$ws.Range("N96").Value = -5627.3333
$ws.Range("H100").Value = 724.9167
$ws.Range("I100").Value = 713
$ws.Range("J100").Value = 748.75
$ws.Range("K100").Value = 1426
$ws.Range("L100").Value = 1497.5
$ws.Range("M100").Value = -885
$ws.Range("N100").Value = -2579.5
$ws.Range("H122").Value = 2910.0588
$ws.Range("J122").Value = 7000
$ws.Range("L122").Value = 21000
$ws.Range("N122").Value = -25900
$ws.Range("H132").Value = 9327.531999999999
$ws.Range("I132").Value = 4990.375
$ws.Range("J132").Value = 18580.133
$ws.Range("K132").Value = 14971.125
$ws.Range("L132").Value = 55740.399
$ws.Range("M132").Value = -12441.125
$ws.Range("N132").Value = -60800.399
